# This script updates the "Pros1-Axl" LR-pairs worksheet (YoungD4) with
# recomputed TPM-based NATMI statistics, per commit "update scripts wuth new tpm".
#
# Underlying ligand (Pros1) "average"/"total" expression values (columns G/H) for each
# Sending cluster and receptor (Axl) "average"/"total" expression values (columns M/N)
# for each Target cluster were recalculated from the refreshed TPM matrix. All of the
# dependent specificity and edge-weight columns (I, J, O, P, Q, R, S, T) are derived
# from those base values:
#   I/J = per-cluster ligand average/total expression, normalized across the 5
#         sending clusters (ligand derived specificity)
#   O/P = per-cluster receptor average/total expression, normalized across the 5
#         target clusters (receptor derived specificity)
#   Q/R = G*M / H*N (edge average/total expression weight)
#   S/T = Q/R normalized across all 25 sending x target cluster combinations
#         (edge derived specificity)
#
# The values below are the refreshed figures taken directly from the updated output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 32.935093
$ws.Range("H2").Value = 98.80527900000001
$ws.Range("I2").Value = 0.247867546118624
$ws.Range("J2").Value = 0.247867546118624
$ws.Range("M2").Value = 3.243447333333334
$ws.Range("N2").Value = 9.730342
$ws.Range("O2").Value = 0.01255381554477167
$ws.Range("P2").Value = 0.01255381554477167
$ws.Range("Q2").Value = 106.8232395639353
$ws.Range("R2").Value = 961.4091560754182
$ws.Range("S2").Value = 0.003111683453508391
$ws.Range("T2").Value = 0.00311168345350839

# Row 3
$ws.Range("G3").Value = 32.935093
$ws.Range("H3").Value = 98.80527900000001
$ws.Range("I3").Value = 0.247867546118624
$ws.Range("J3").Value = 0.247867546118624
$ws.Range("O3").Value = 0.2470858318854151
$ws.Range("P3").Value = 0.2470858318854151
$ws.Range("Q3").Value = 2102.508908006265
$ws.Range("R3").Value = 18922.58017205638
$ws.Range("S3").Value = 0.06124455883011671
$ws.Range("T3").Value = 0.0612445588301167

# Row 4
$ws.Range("G4").Value = 32.935093
$ws.Range("H4").Value = 98.80527900000001
$ws.Range("I4").Value = 0.247867546118624
$ws.Range("J4").Value = 0.247867546118624
$ws.Range("M4").Value = 103.7552793333333
$ws.Range("N4").Value = 311.265838
$ws.Range("O4").Value = 0.4015864926064037
$ws.Range("P4").Value = 0.4015864926064037
$ws.Range("Q4").Value = 3417.189774084311
$ws.Range("R4").Value = 30754.7079667588
$ws.Range("S4").Value = 0.09954025847673423
$ws.Range("T4").Value = 0.09954025847673421

# Row 5
$ws.Range("G5").Value = 32.935093
$ws.Range("H5").Value = 98.80527900000001
$ws.Range("I5").Value = 0.247867546118624
$ws.Range("J5").Value = 0.247867546118624
$ws.Range("M5").Value = 23.61337433333334
$ws.Range("N5").Value = 70.840123
$ws.Range("O5").Value = 0.09139594860190291
$ws.Range("P5").Value = 0.0913959486019029
$ws.Range("Q5").Value = 777.7086797121465
$ws.Range("R5").Value = 6999.378117409318
$ws.Range("S5").Value = 0.02265408950513756
$ws.Range("T5").Value = 0.02265408950513755

# Row 6
$ws.Range("G6").Value = 32.935093
$ws.Range("H6").Value = 98.80527900000001
$ws.Range("I6").Value = 0.247867546118624
$ws.Range("J6").Value = 0.247867546118624
$ws.Range("M6").Value = 63.91341533333334
$ws.Range("N6").Value = 191.740246
$ws.Range("O6").Value = 0.2473779113615065
$ws.Range("P6").Value = 0.2473779113615065
$ws.Range("Q6").Value = 2104.99427795096
$ws.Range("R6").Value = 18944.94850155864
$ws.Range("S6").Value = 0.0613169558531271
$ws.Range("T6").Value = 0.06131695585312708

# Row 7
$ws.Range("I7").Value = 0.4962147730988433
$ws.Range("J7").Value = 0.4962147730988432
$ws.Range("M7").Value = 3.243447333333334
$ws.Range("N7").Value = 9.730342
$ws.Range("O7").Value = 0.01255381554477167
$ws.Range("P7").Value = 0.01255381554477167
$ws.Range("Q7").Value = 213.8532067305553
$ws.Range("R7").Value = 1924.678860574998
$ws.Range("S7").Value = 0.006229388732073605
$ws.Range("T7").Value = 0.006229388732073604

# Row 8
$ws.Range("I8").Value = 0.4962147730988433
$ws.Range("J8").Value = 0.4962147730988432
$ws.Range("O8").Value = 0.2470858318854151
$ws.Range("P8").Value = 0.2470858318854151
$ws.Range("S8").Value = 0.1226076400049602
$ws.Range("T8").Value = 0.1226076400049602

# Row 9
$ws.Range("I9").Value = 0.4962147730988433
$ws.Range("J9").Value = 0.4962147730988432
$ws.Range("M9").Value = 103.7552793333333
$ws.Range("N9").Value = 311.265838
$ws.Range("O9").Value = 0.4015864926064037
$ws.Range("P9").Value = 0.4015864926064037
$ws.Range("Q9").Value = 6840.99259840749
$ws.Range("R9").Value = 61568.93338566741
$ws.Range("S9").Value = 0.1992731503082469
$ws.Range("T9").Value = 0.1992731503082469

# Row 10
$ws.Range("I10").Value = 0.4962147730988433
$ws.Range("J10").Value = 0.4962147730988432
$ws.Range("M10").Value = 23.61337433333334
$ws.Range("N10").Value = 70.840123
$ws.Range("O10").Value = 0.09139594860190291
$ws.Range("P10").Value = 0.0913959486019029
$ws.Range("Q10").Value = 1556.922405064176
$ws.Range("R10").Value = 14012.30164557759
$ws.Range("S10").Value = 0.04535201989764679
$ws.Range("T10").Value = 0.04535201989764678

# Row 11
$ws.Range("I11").Value = 0.4962147730988433
$ws.Range("J11").Value = 0.4962147730988432
$ws.Range("M11").Value = 63.91341533333334
$ws.Range("N11").Value = 191.740246
$ws.Range("O11").Value = 0.2473779113615065
$ws.Range("P11").Value = 0.2473779113615065
$ws.Range("Q11").Value = 4214.06220525502
$ws.Range("R11").Value = 37926.55984729518
$ws.Range("S11").Value = 0.1227525741559157
$ws.Range("T11").Value = 0.1227525741559157

# Row 12
$ws.Range("G12").Value = 13.46314666666667
$ws.Range("H12").Value = 40.38944
$ws.Range("I12").Value = 0.1013228390550407
$ws.Range("J12").Value = 0.1013228390550407
$ws.Range("M12").Value = 3.243447333333334
$ws.Range("N12").Value = 9.730342
$ws.Range("O12").Value = 0.01255381554477167
$ws.Range("P12").Value = 0.01255381554477167
$ws.Range("Q12").Value = 43.66700715427556
$ws.Range("R12").Value = 393.00306438848
$ws.Range("S12").Value = 0.001271988231969568
$ws.Range("T12").Value = 0.001271988231969568

# Row 13
$ws.Range("G13").Value = 13.46314666666667
$ws.Range("H13").Value = 40.38944
$ws.Range("I13").Value = 0.1013228390550407
$ws.Range("J13").Value = 0.1013228390550407
$ws.Range("O13").Value = 0.2470858318854151
$ws.Range("P13").Value = 0.2470858318854151
$ws.Range("Q13").Value = 859.4597196510578
$ws.Range("R13").Value = 7735.13747685952
$ws.Range("S13").Value = 0.02503543797690677
$ws.Range("T13").Value = 0.02503543797690677

# Row 14
$ws.Range("G14").Value = 13.46314666666667
$ws.Range("H14").Value = 40.38944
$ws.Range("I14").Value = 0.1013228390550407
$ws.Range("J14").Value = 0.1013228390550407
$ws.Range("M14").Value = 103.7552793333333
$ws.Range("N14").Value = 311.265838
$ws.Range("O14").Value = 0.4015864926064037
$ws.Range("P14").Value = 0.4015864926064037
$ws.Range("Q14").Value = 1396.872543105635
$ws.Range("R14").Value = 12571.85288795072
$ws.Range("S14").Value = 0.04068988355703695
$ws.Range("T14").Value = 0.04068988355703695

# Row 15
$ws.Range("G15").Value = 13.46314666666667
$ws.Range("H15").Value = 40.38944
$ws.Range("I15").Value = 0.1013228390550407
$ws.Range("J15").Value = 0.1013228390550407
$ws.Range("M15").Value = 23.61337433333334
$ws.Range("N15").Value = 70.840123
$ws.Range("O15").Value = 0.09139594860190291
$ws.Range("P15").Value = 0.0913959486019029
$ws.Range("Q15").Value = 317.9103219445689
$ws.Range("R15").Value = 2861.19289750112
$ws.Range("S15").Value = 0.009260496990473384
$ws.Range("T15").Value = 0.009260496990473382

# Row 16
$ws.Range("G16").Value = 13.46314666666667
$ws.Range("H16").Value = 40.38944
$ws.Range("I16").Value = 0.1013228390550407
$ws.Range("J16").Value = 0.1013228390550407
$ws.Range("M16").Value = 63.91341533333334
$ws.Range("N16").Value = 191.740246
$ws.Range("O16").Value = 0.2473779113615065
$ws.Range("P16").Value = 0.2473779113615065
$ws.Range("Q16").Value = 860.475684600249
$ws.Range("R16").Value = 7744.28116140224
$ws.Range("S16").Value = 0.02506503229865406
$ws.Range("T16").Value = 0.02506503229865405

# Row 17
$ws.Range("G17").Value = 6.098311666666667
$ws.Range("H17").Value = 18.294935
$ws.Range("I17").Value = 0.04589552998326869
$ws.Range("J17").Value = 0.04589552998326869
$ws.Range("M17").Value = 3.243447333333334
$ws.Range("N17").Value = 9.730342
$ws.Range("O17").Value = 0.01255381554477167
$ws.Range("P17").Value = 0.01255381554477167
$ws.Range("Q17").Value = 19.77955271308556
$ws.Range("R17").Value = 178.01597441777
$ws.Range("S17").Value = 0.0005761640177394928
$ws.Range("T17").Value = 0.0005761640177394927

# Row 18
$ws.Range("G18").Value = 6.098311666666667
$ws.Range("H18").Value = 18.294935
$ws.Range("I18").Value = 0.04589552998326869
$ws.Range("J18").Value = 0.04589552998326869
$ws.Range("O18").Value = 0.2470858318854151
$ws.Range("P18").Value = 0.2470858318854151
$ws.Range("Q18").Value = 389.3037315232478
$ws.Range("R18").Value = 3503.733583709231
$ws.Range("S18").Value = 0.01134013520573796
$ws.Range("T18").Value = 0.01134013520573796

# Row 19
$ws.Range("G19").Value = 6.098311666666667
$ws.Range("H19").Value = 18.294935
$ws.Range("I19").Value = 0.04589552998326869
$ws.Range("J19").Value = 0.04589552998326869
$ws.Range("M19").Value = 103.7552793333333
$ws.Range("N19").Value = 311.265838
$ws.Range("O19").Value = 0.4015864926064037
$ws.Range("P19").Value = 0.4015864926064037
$ws.Range("Q19").Value = 632.7320304367256
$ws.Range("R19").Value = 5694.588273930531
$ws.Range("S19").Value = 0.01843102491229291
$ws.Range("T19").Value = 0.01843102491229291

# Row 20
$ws.Range("G20").Value = 6.098311666666667
$ws.Range("H20").Value = 18.294935
$ws.Range("I20").Value = 0.04589552998326869
$ws.Range("J20").Value = 0.04589552998326869
$ws.Range("M20").Value = 23.61337433333334
$ws.Range("N20").Value = 70.840123
$ws.Range("O20").Value = 0.09139594860190291
$ws.Range("P20").Value = 0.0913959486019029
$ws.Range("Q20").Value = 144.0017161863339
$ws.Range("R20").Value = 1296.015445677005
$ws.Range("S20").Value = 0.004194665499407919
$ws.Range("T20").Value = 0.004194665499407919

# Row 21
$ws.Range("G21").Value = 6.098311666666667
$ws.Range("H21").Value = 18.294935
$ws.Range("I21").Value = 0.04589552998326869
$ws.Range("J21").Value = 0.04589552998326869
$ws.Range("M21").Value = 63.91341533333334
$ws.Range("N21").Value = 191.740246
$ws.Range("O21").Value = 0.2473779113615065
$ws.Range("P21").Value = 0.2473779113615065
$ws.Range("Q21").Value = 389.763926383779
$ws.Range("R21").Value = 3507.875337454011
$ws.Range("S21").Value = 0.01135354034809041
$ws.Range("T21").Value = 0.01135354034809041

# Row 22
$ws.Range("G22").Value = 14.44328633333333
$ws.Range("H22").Value = 43.329859
$ws.Range("I22").Value = 0.1086993117442235
$ws.Range("J22").Value = 0.1086993117442234
$ws.Range("M22").Value = 3.243447333333334
$ws.Range("N22").Value = 9.730342
$ws.Range("O22").Value = 0.01255381554477167
$ws.Range("P22").Value = 0.01255381554477167
$ws.Range("Q22").Value = 46.84603854241978
$ws.Range("R22").Value = 421.614346881778
$ws.Range("S22").Value = 0.001364591109480614
$ws.Range("T22").Value = 0.001364591109480614

# Row 23
$ws.Range("G23").Value = 14.44328633333333
$ws.Range("H23").Value = 43.329859
$ws.Range("I23").Value = 0.1086993117442235
$ws.Range("J23").Value = 0.1086993117442234
$ws.Range("O23").Value = 0.2470858318854151
$ws.Range("P23").Value = 0.2470858318854151
$ws.Range("Q23").Value = 922.0298292984469
$ws.Range("R23").Value = 8298.268463686021
$ws.Range("S23").Value = 0.02685805986769353
$ws.Range("T23").Value = 0.02685805986769352

# Row 24
$ws.Range("G24").Value = 14.44328633333333
$ws.Range("H24").Value = 43.329859
$ws.Range("I24").Value = 0.1086993117442235
$ws.Range("J24").Value = 0.1086993117442234
$ws.Range("M24").Value = 103.7552793333333
$ws.Range("N24").Value = 311.265838
$ws.Range("O24").Value = 0.4015864926064037
$ws.Range("P24").Value = 0.4015864926064037
$ws.Range("Q24").Value = 1498.567208006315
$ws.Range("R24").Value = 13487.10487205684
$ws.Range("S24").Value = 0.04365217535209277
$ws.Range("T24").Value = 0.04365217535209275

# Row 25
$ws.Range("G25").Value = 14.44328633333333
$ws.Range("H25").Value = 43.329859
$ws.Range("I25").Value = 0.1086993117442235
$ws.Range("J25").Value = 0.1086993117442234
$ws.Range("M25").Value = 23.61337433333334
$ws.Range("N25").Value = 70.840123
$ws.Range("O25").Value = 0.09139594860190291
$ws.Range("P25").Value = 0.0913959486019029
$ws.Range("Q25").Value = 341.0547267925175
$ws.Range("R25").Value = 3069.492541132657
$ws.Range("S25").Value = 0.009934676709237268
$ws.Range("T25").Value = 0.009934676709237264

# Row 26
$ws.Range("G26").Value = 14.44328633333333
$ws.Range("H26").Value = 43.329859
$ws.Range("I26").Value = 0.1086993117442235
$ws.Range("J26").Value = 0.1086993117442234
$ws.Range("M26").Value = 63.91341533333334
$ws.Range("N26").Value = 191.740246
$ws.Range("O26").Value = 0.2473779113615065
$ws.Range("P26").Value = 0.2473779113615065
$ws.Range("Q26").Value = 923.1197582005905
$ws.Range("R26").Value = 8308.077823805314
$ws.Range("S26").Value = 0.02688980870571928
$ws.Range("T26").Value = 0.02688980870571927

